$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Final Results")

# Make "Final Results" the active sheet/tab (was "IPC-anagram" chartsheet).
$ws.Activate()

# Set initial IPC values for the "IPC" and "Init" graphs (was interim/carried
# values; now a flat initial baseline of 100).
$ws.Range("B38").Value = 100

$ws.Range("B39").ClearFormats()
$ws.Range("B39").Value = 100

$ws.Range("B40").ClearFormats()
$ws.Range("B40").Value = 100

$ws.Range("B41").ClearFormats()
$ws.Range("B41").Value = 100

$ws.Range("B44").Value = 100
$ws.Range("B45").Value = 100
$ws.Range("B46").Value = 100
$ws.Range("B47").Value = 100

# Move the on-sheet selection to B50 (was C39).
$ws.Range("B50").Select() | Out-Null
